$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: num_customers 49 -> 50, retention_rate recalculated (50/2252)
$ws.Range("C27").Value = 50
$ws.Range("E27").Value = 0.02220248667850799

# Row 31: num_customers 51 -> 53, retention_rate recalculated (53/2312)
$ws.Range("C31").Value = 53
$ws.Range("E31").Value = 0.02292387543252595

# Row 34: num_customers 85 -> 86, retention_rate recalculated (86/2256)
$ws.Range("C34").Value = 86
$ws.Range("E34").Value = 0.03812056737588652

# Row 37: num_customers and cohort_size 904 -> 913 (retention_rate stays 1)
$ws.Range("C37").Value = 913
$ws.Range("D37").Value = 913
